$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 3 / Row 8 header cells: "time worked" / "total time" / "hours" labels
# These cells previously just carried the date-row style (s="1") with no
# content; they now hold text and lose that formatting.
# ---------------------------------------------------------------------------
$ws.Range("D3").ClearFormats()
$ws.Range("D3").Value = "time worked"
$ws.Range("F3").ClearFormats()
$ws.Range("F3").Value = "total time"
$ws.Range("G3").ClearFormats()
$ws.Range("G3").Value = "hours"

# ---------------------------------------------------------------------------
# Row 4 / Row 5 : numeric minutes for the first day's two entries
# ---------------------------------------------------------------------------
$ws.Range("D4").Value = 30
$ws.Range("F4").Value = 150
$ws.Range("G4").Value = "2,5"

$ws.Range("D5").Value = 120

# ---------------------------------------------------------------------------
# Row 8 : second day header labels
# ---------------------------------------------------------------------------
$ws.Range("D8").Value = "time worked"
$ws.Range("F8").Value = "total time"
$ws.Range("G8").Value = "hours"

# Row 9 : second day entry totals
$ws.Range("D9").Value = 120
$ws.Range("F9").Value = 120
$ws.Range("G9").Value = 2

# ---------------------------------------------------------------------------
# Row 11 : new day (26 Dec 2017) header row
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = 43095
$ws.Range("B11").NumberFormat = "d-mmm"
$ws.Range("D11").Value = "time worked"
$ws.Range("F11").Value = "total time"
$ws.Range("G11").Value = "hours"

# Row 12 : movement navmesh
$ws.Range("B12").Value = "13:35-14:07"
$ws.Range("C12").Value = "movement navmesh"
$ws.Range("D12").Value = 32
$ws.Range("F12").Formula = "=D12+D13+D14+D15+D16+D17"
$ws.Range("G12").Formula = "=F12/60"

# Row 13 : camera
$ws.Range("B13").Value = "14:33-15:00"
$ws.Range("C13").Value = "camera"
$ws.Range("D13").Value = 27

# Row 14 : player switching
$ws.Range("B14").Value = "14:15-14:32"
$ws.Range("C14").Value = "player switching"
$ws.Range("D14").Value = 17

# Row 15 : camera + player switching
$ws.Range("B15").Value = "16:20-17:35"
$ws.Range("C15").Value = "camera + player switching"
$ws.Range("D15").Value = 75

# Row 16 : HUD
$ws.Range("B16").Value = "19:30-21:15"
$ws.Range("C16").Value = "HUD"
$ws.Range("D16").Value = 105

# Row 17 : stats
$ws.Range("B17").Value = "21:30-22:30"
$ws.Range("C17").Value = "stats"
$ws.Range("D17").Value = 60

# ---------------------------------------------------------------------------
# Row 19 : new day (27 Dec 2017) header row
# ---------------------------------------------------------------------------
$ws.Range("B19").Value = 43096
$ws.Range("B19").NumberFormat = "d-mmm"
$ws.Range("D19").Value = "time worked"
$ws.Range("F19").Value = "total time"
$ws.Range("G19").Value = "hours"

# Row 20 : weapons
$ws.Range("B20").Value = "13:30-13:40"
$ws.Range("C20").Value = "weapons"
$ws.Range("D20").Value = 10
$ws.Range("F20").Formula = "=D20+D21+D22+D24+D23+D25+D26"
$ws.Range("G20").Formula = "=F20/60"

# Row 21 : inventory
$ws.Range("B21").Value = "13:40-15:44"
$ws.Range("C21").Value = "inventory"
$ws.Range("D21").Value = 124

# Row 22 : enemy
$ws.Range("B22").Value = "15:50-16:22"
$ws.Range("C22").Value = "enemy"
$ws.Range("D22").Value = 32

# Row 23 : enemy attack + damage system on the player
$ws.Range("B23").Value = "17:20-18:00"
$ws.Range("C23").Value = "enemy attack + damage system on the player"
$ws.Range("D23").Value = 40

# Row 24 : enemy attack + damage system on the player
$ws.Range("B24").Value = "19:34-19:54"
$ws.Range("C24").Value = "enemy attack + damage system on the player"
$ws.Range("D24").Value = 20

# Row 25 : weapons
$ws.Range("B25").Value = "20:00-20:38"
$ws.Range("C25").Value = "weapons"
$ws.Range("D25").Value = 38

# Row 26 : weapon switching + player stats + inventory
$ws.Range("B26").Value = "21:30-00:30"
$ws.Range("C26").Value = "weapon switching + player stats + inventory"
$ws.Range("D26").Value = 180

# ---------------------------------------------------------------------------
# Row 28 : new day (28 Dec 2017) date marker
# ---------------------------------------------------------------------------
$ws.Range("B28").Value = 43097
$ws.Range("B28").NumberFormat = "d-mmm"

# Extra column widths used by the new D/E/F columns of data
$ws.Columns.Item(4).ColumnWidth = 16.666666666666668
$ws.Columns.Item(5).ColumnWidth = 9.666666666666666
$ws.Columns.Item(6).ColumnWidth = 12.0

# Re-create the scroll position / active selection recorded for the sheet
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("J21").Select()
